$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM (kidney)")
$ws.Activate()

$ws.Range("A8").Value = "Rodrigues et al., 1983"
$ws.Range("B8").Value = 91
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 12
$ws.Range("E8").Formula = "=C8/SQRT(D8)"

$table = $ws.ListObjects.Item("Table58")
$table.Resize($ws.Range("A1:E8"))

$ws.Range("E8").Select()
